# Connected Office Test Data - updated device testing procedures
#
# 1. On the "Test Results" sheet, mark every CRUD test result
#    (Create/Read/Update/Delete, columns B:E, rows 2-24) as Passed (TRUE).
# 2. Update the remembered cell selections on the "Device" and
#    "Test Results" sheets to reflect where the user left off,
#    leaving "Test Results" as the active sheet.

$wb = $excel.ActiveWorkbook

# --- Device sheet: move remembered selection from C38 to F2 -----------
$wsDevice = $wb.Worksheets.Item("Device")
$wsDevice.Activate()
$wsDevice.Range("F2").Select()

# --- Test Results sheet: flip all CRUD test results to Passed ---------
$wsResults = $wb.Worksheets.Item("Test Results")
$wsResults.Activate()
$wsResults.Range("B2:E24").Value = $true

# move remembered selection from E2 to K13, keep this sheet active
$wsResults.Range("K13").Select()
